$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Fill-EmptyCell($row, $col) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $rng.Text = "S"
    $rng2 = $cell.Range
    $rng2.Font.Size = 8
    $rng2.Font.SizeBi = 8
}

function Replace-CellChar($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $charRng = $cell.Range.Characters.Item(1)
    $charRng.Text = $newText
}

# Row "2" (table row 3): CancelButton, PowerButton, TimeButton, Door were empty -> "S"
Fill-EmptyCell 3 2
Fill-EmptyCell 3 3
Fill-EmptyCell 3 4
Fill-EmptyCell 3 5

# Row "3" (table row 4): CancelButton, PowerButton, TimeButton, Door were empty -> "S"
Fill-EmptyCell 4 2
Fill-EmptyCell 4 3
Fill-EmptyCell 4 4
Fill-EmptyCell 4 5

# Row "4*" (table row 5), Output column: "X" -> "S"
Replace-CellChar 5 12 "S"

# Row "5*" (table row 6), Output column: "X" -> "S"
Replace-CellChar 6 12 "S"

# Row "6*" (table row 7), Output column: "X" -> "S"
Replace-CellChar 7 12 "S"
